$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 90: Katie Harding / Freenome / https://www.freenome.com/
$ws.Range("A90").Value = "Katie Harding"
$ws.Range("B90").Value = "Freenome"
$ws.Range("C90").Value = "https://www.freenome.com/ "

# Row 91: Youn Kyeong Chang / FDA / https://www.fda.gov/
$ws.Range("A91").Value = "Youn Kyeong Chang"
$ws.Range("B91").Value = "FDA"
$ws.Range("C91").Value = "https://www.fda.gov/"
